# Chapter 1 vocabulary list - finish Ch1 (11-19) ✨
# Adds the remaining word_EN/word_AR pairs + page numbers for pages 9-19,
# re-applies the sheet's AutoFilter header and widens column B to fit the
# longest Arabic translation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add missing page numbers (column C) to existing rows 26-31, 33 ---
$ws.Range("C26").Value = 9
$ws.Range("C27").Value = 9
$ws.Range("C28").Value = 9
$ws.Range("C29").Value = 9
$ws.Range("C30").Value = 9
$ws.Range("C31").Value = 9
$ws.Range("C33").Value = 10

# --- Row 34: fill in A/B (C34/D34 already existed) ---
$ws.Range("A34").Value = "despite"
$ws.Range("B34").Value = "بالرغم من"

# --- New rows 35-64: word_EN / word_AR / page ---
$ws.Range("A35").Value = "acronym"
$ws.Range("B35").Value = "اختصار"
$ws.Range("C35").Value = 11
$ws.Range("A36").Value = "scares"
$ws.Range("B36").Value = "يخيف"
$ws.Range("C36").Value = 11
$ws.Range("A37").Value = "significance"
$ws.Range("B37").Value = "دلالة"
$ws.Range("C37").Value = 12
$ws.Range("A38").Value = "stochastic"
$ws.Range("B38").Value = "العشوائية"
$ws.Range("C38").Value = 12
$ws.Range("A39").Value = "instruct"
$ws.Range("B39").Value = "إرشاد"
$ws.Range("C39").Value = 13
$ws.Range("A40").Value = "fancy"
$ws.Range("B40").Value = "خيالي؟"
$ws.Range("C40").Value = 13
$ws.Range("A41").Value = "exhibited"
$ws.Range("B41").Value = "عرضت"
$ws.Range("C41").Value = 13
$ws.Range("A42").Value = "obstinacy"
$ws.Range("B42").Value = "عناد"
$ws.Range("C42").Value = 13
$ws.Range("A43").Value = "sympathize"
$ws.Range("B43").Value = "تتعاطف"
$ws.Range("C43").Value = 13
$ws.Range("A44").Value = "comprise"
$ws.Range("B44").Value = "تضم"
$ws.Range("C44").Value = 13
$ws.Range("A45").Value = "extent"
$ws.Range("B45").Value = "مدى"
$ws.Range("C45").Value = 13
$ws.Range("A46").Value = "expertise"
$ws.Range("B46").Value = "خبرة"
$ws.Range("C46").Value = 13
$ws.Range("A47").Value = "bolster"
$ws.Range("B47").Value = "دعم"
$ws.Range("C47").Value = 14
$ws.Range("A48").Value = "narrative"
$ws.Range("B48").Value = "رواية"
$ws.Range("C48").Value = 15
$ws.Range("A49").Value = "sophisticated"
$ws.Range("B49").Value = "متطور"
$ws.Range("C49").Value = 15
$ws.Range("A50").Value = "disempowering"
$ws.Range("B50").Value = "إضعاف"
$ws.Range("C50").Value = 15
$ws.Range("A51").Value = "meticulous"
$ws.Range("B51").Value = "شديد الإنتباه للتفاصيل"
$ws.Range("C51").Value = 16
$ws.Range("A52").Value = "groundbreaking"
$ws.Range("B52").Value = "رائدة"
$ws.Range("C52").Value = 16
$ws.Range("A53").Value = "unparalleled"
$ws.Range("B53").Value = "لا مثيل لها"
$ws.Range("C53").Value = 16
$ws.Range("A54").Value = "inquisitiveness"
$ws.Range("B54").Value = "الفضول"
$ws.Range("C54").Value = 16
$ws.Range("A55").Value = "messy data"
$ws.Range("B55").Value = "بيانات فوضوية"
$ws.Range("C55").Value = 16
$ws.Range("A56").Value = "sheer perseverance"
$ws.Range("B56").Value = "المثابرة المطلقة"
$ws.Range("C56").Value = 16
$ws.Range("A57").Value = "gravitate"
$ws.Range("B57").Value = "تنجذب"
$ws.Range("C57").Value = 16
$ws.Range("A58").Value = "gifted"
$ws.Range("B58").Value = "موهوبين"
$ws.Range("C58").Value = 16
$ws.Range("A59").Value = "savvy"
$ws.Range("B59").Value = "الدهاء"
$ws.Range("C59").Value = 16
$ws.Range("A60").Value = "acumen"
$ws.Range("B60").Value = "الفطنة"
$ws.Range("C60").Value = 16
$ws.Range("A61").Value = "evolved"
$ws.Range("B61").Value = "تطورت"
$ws.Range("C61").Value = 16
$ws.Range("A62").Value = "entrepreneur"
$ws.Range("B62").Value = "رائد أعمال"
$ws.Range("C62").Value = 17
$ws.Range("A63").Value = "risk-tolerant"
$ws.Range("B63").Value = "متسامح مع المخاطر"
$ws.Range("C63").Value = 17
$ws.Range("A64").Value = "expertise"
$ws.Range("B64").Value = "خبرة"
$ws.Range("C64").Value = 17

# --- Row 65: page 19, date 26-4-2023 (reuses existing date string) ---
$ws.Range("C65").Value = 19
$ws.Range("D65").Value = "26-4-2023"

# --- Formatting / view touch-ups ---

# Column B needs to be a bit wider to fit the new (longer) Arabic text
$ws.Columns.Item(2).ColumnWidth = 14.3

# Re-apply the header AutoFilter (A1:D1) and its hidden _FilterDatabase name
[void]$ws.Range("A1:D1").AutoFilter()
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$D`$1")
$filterName.Visible = $false

# Scroll the view down to where we left off and select the last edited cell
$win = $excel.ActiveWindow
$win.ScrollRow = 40
$win.ScrollColumn = 1
[void]$excel.Goto($ws.Range("A40"), $true)
[void]$ws.Range("D62").Select()
